$wb = $excel.ActiveWorkbook

# The workbook currently has 3 sheets: 2021-Q2, 2021-Q4, 总计 (Total).
# We need to end up with 4 sheets: 2021-Q2, 2021-Q4, 2022-Q1, 总计
#   - "2022-Q1" is a brand-new sheet with fund-holding detail rows, inserted
#     in the slot formerly occupied by "总计" (reusing its sheetId/rId).
#   - "总计" becomes a new physical sheet (after "2022-Q1") that keeps the
#     same summary-table layout as before, but with a new row inserted at
#     the top for the "2022-Q1" totals (pushing the older rows down).

$wsOld = $wb.Worksheets.Item(3)

# --- Step 1: clone the existing "总计" sheet BEFORE we touch it, so the
# clone keeps its exact data/formatting; this clone will become the new
# "总计" sheet. It is placed immediately after $wsOld.
$wsOld.Copy($null, $wsOld)
$wsTotal = $wb.Worksheets.Item(4)
$wsTotal.Name = "总计NEW"

# --- Step 2: turn the original sheet into "2022-Q1" and replace its
# contents with the fund holding data.
$wsOld.Name = "2022-Q1"

# Extend the header style (s="2", already on B1:D1) across the new columns.
$wsOld.Range("D1").Copy()
$wsOld.Range("E1:H1").PasteSpecial(-4122)

$wsOld.Range("B1").Value = "基金代码"
$wsOld.Range("C1").Value = "基金名称"
$wsOld.Range("D1").Value = "基金规模"
$wsOld.Range("E1").Value = "股票总仓位"
$wsOld.Range("F1").Value = "仓位占比"
$wsOld.Range("G1").Value = "持有市值(亿元)"
$wsOld.Range("H1").Value = "仓位排名"

# Extend the index-column style (s="2", already on A2:A3) down to A4:A5.
$wsOld.Range("A2").Copy()
$wsOld.Range("A4:A5").PasteSpecial(-4122)

$wsOld.Range("A2").Value = 0
$wsOld.Range("B2").Value = "'003713"
$wsOld.Range("B2").Style = "Normal"
$wsOld.Range("C2").Value = "英大睿盛灵活配置混合A"
$wsOld.Range("D2").Value = "'5.99"
$wsOld.Range("D2").Style = "Normal"
$wsOld.Range("E2").Value = "'87.42"
$wsOld.Range("E2").Style = "Normal"
$wsOld.Range("F2").Value = "'7.52"
$wsOld.Range("F2").Style = "Normal"
$wsOld.Range("G2").Value = "'0.4504"
$wsOld.Range("G2").Style = "Normal"
$wsOld.Range("H2").Value = 1

$wsOld.Range("A3").Value = 1
$wsOld.Range("B3").Value = "'003714"
$wsOld.Range("B3").Style = "Normal"
$wsOld.Range("C3").Value = "英大睿盛灵活配置混合C"
$wsOld.Range("D3").Value = "'2.40"
$wsOld.Range("D3").Style = "Normal"
$wsOld.Range("E3").Value = "'87.42"
$wsOld.Range("E3").Style = "Normal"
$wsOld.Range("F3").Value = "'7.52"
$wsOld.Range("F3").Style = "Normal"
$wsOld.Range("G3").Value = "'0.1805"
$wsOld.Range("G3").Style = "Normal"
$wsOld.Range("H3").Value = 1

$wsOld.Range("A4").Value = 2
$wsOld.Range("B4").Value = "'003446"
$wsOld.Range("B4").Style = "Normal"
$wsOld.Range("C4").Value = "英大睿鑫灵活配置混合A"
$wsOld.Range("D4").Value = "'0.59"
$wsOld.Range("D4").Style = "Normal"
$wsOld.Range("E4").Value = "'89.46"
$wsOld.Range("E4").Style = "Normal"
$wsOld.Range("F4").Value = "'6.15"
$wsOld.Range("F4").Style = "Normal"
$wsOld.Range("G4").Value = "'0.0363"
$wsOld.Range("G4").Style = "Normal"
$wsOld.Range("H4").Value = 5

$wsOld.Range("A5").Value = 3
$wsOld.Range("B5").Value = "'003447"
$wsOld.Range("B5").Style = "Normal"
$wsOld.Range("C5").Value = "英大睿鑫灵活配置混合C"
$wsOld.Range("D5").Value = "'0.51"
$wsOld.Range("D5").Style = "Normal"
$wsOld.Range("E5").Value = "'89.46"
$wsOld.Range("E5").Style = "Normal"
$wsOld.Range("F5").Value = "'6.15"
$wsOld.Range("F5").Style = "Normal"
$wsOld.Range("G5").Value = "'0.0314"
$wsOld.Range("G5").Style = "Normal"
$wsOld.Range("H5").Value = 5

# --- Step 3: finish the new "总计" sheet: insert a fresh row 2 for the
# 2022-Q1 totals and renumber/rewrite the (now shifted) old rows so every
# cell gets an explicit, known-good value.
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("B2:D2").Style = "Normal"

$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 4
$wsTotal.Range("D2").Value = 0.7

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2021-Q4"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.34

$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2021-Q2"
$wsTotal.Range("C4").Value = 1
$wsTotal.Range("D4").Value = 0.01

$wsTotal.Name = "总计"

# Keep the original active tab (2021-Q2) selected, matching the source file.
$wb.Worksheets.Item(1).Activate()
